# Append new user rows (12-25) to the "users" sheet.
# Column D holds phone numbers that start with a literal "+" sign (and, in
# one case, a 19-digit number) which Excel would otherwise auto-convert to
# a plain number (stripping the "+", and/or losing precision). To keep
# these as text we temporarily force a Text ("@") number format before
# assigning the value, then clear the formatting again afterwards so the
# cell is left with the default style (matching freshly authored rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-08-07 12:31:33", "boris",     "Benin",       "+2294545455"),
    @("2024-08-07 12:42:41", "boris",     "Bolivia",     "+5914545455"),
    @("2024-08-07 14:30:47", "pai natal", "Afghanistan", "+93777"),
    @("2024-08-07 14:34:02", "ddd",       "Bolivia",     "+591999"),
    @("2024-08-07 14:45:55", "ddd",       "Bolivia",     "+5919999"),
    @("2024-08-07 14:49:23", "ddd",       "Bolivia",     "+59199991"),
    @("2024-08-07 14:58:30", "qqq",       "Afghanistan", "+93001100"),
    @("2024-08-07 15:00:06", "qqq",       "Afghanistan", "+930011001"),
    @("2024-08-07 15:00:20", "qqq",       "Afghanistan", "+9300110011"),
    @("2024-08-07 15:00:42", "qqq",       "Afghanistan", "+93001100110"),
    @("2024-08-07 15:02:55", "aa",        "Afghanistan", "+931"),
    @("2024-08-07 15:03:25", "aa",        "Afghanistan", "+9310"),
    @("2024-08-07 16:41:29", "jhg",       "Afghanistan", "+93887788"),
    @("2024-08-07 16:45:02", "boris",     "Afghanistan", "+933211654987654321")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Force the phone-number cell to Text first so the leading "+" (and
    # full digit precision) survives, then strip the formatting again.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[3]
    $dCell.ClearFormats()
}

Write-Host "Added rows 12-25"
